$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.536.73'
$ws.Range('E2').Value = '  -1.24%  '

$ws.Range('D3').Value = '2.653.52'
$ws.Range('E3').Value = '  +1.35%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '''518.36'
$ws.Range('E5').Value = '  -0.28%  '

$ws.Range('D6').Value = '''149.38'
$ws.Range('E6').Value = '  -0.69%  '

$ws.Range('D7').Value = '''0.994'
$ws.Range('E7').Value = '  -0.51%  '

$ws.Range('D8').Value = '''0.577'
$ws.Range('E8').Value = '  +0.50%  '

$ws.Range('D9').Value = '2.691.06'
$ws.Range('E9').Value = '  +2.65%  '

$ws.Range('D10').Value = '''6.59'
$ws.Range('E10').Value = '  +3.25%  '

$ws.Range('E11').Value = '  +0.46%  '

$ws.Range('D12').Value = '''0.342'
$ws.Range('E12').Value = '  -0.39%  '

$ws.Range('E13').Value = '  -1.07%  '

$ws.Range('D14').Value = '3.112.21'
$ws.Range('E14').Value = '  +0.91%  '

$ws.Range('D15').Value = '59.363.51'
$ws.Range('E15').Value = '  -1.55%  '

$ws.Range('D16').Value = '''21.57'
$ws.Range('E16').Value = '  +0.47%  '

$ws.Range('E17').Value = '  +0.89%  '

$ws.Range('D18').Value = '2.673.59'
$ws.Range('E18').Value = '  +1.47%  '

$ws.Range('D19').Value = '''4.65'
$ws.Range('E19').Value = '  -0.21%  '

$ws.Range('D20').Value = '''347.93'
$ws.Range('E20').Value = '  +0.39%  '

$ws.Range('D21').Value = '''10.60'
$ws.Range('E21').Value = '  +1.86%  '

$ws.Range('D22').Value = '''6.23'
$ws.Range('E22').Value = '  +1.22%  '

$ws.Range('E23').Value = '  +0.44%  '

$ws.Range('D24').Value = '''61.07'
$ws.Range('E24').Value = '  +0.24%  '

$ws.Range('D25').Value = '''0.429'
$ws.Range('E25').Value = '  +1.89%  '

$ws.Range('D26').Value = '2.764.62'
$ws.Range('E26').Value = '  +0.79%  '

$ws.Range('D27').Value = '''0.991'
$ws.Range('E27').Value = '  -0.85%  '

$ws.Range('E28').Value = '  -1.67%  '

$ws.Range('E29').Value = '  +1.29%  '

$ws.Range('D30').Value = '''7.20'
$ws.Range('E30').Value = '  +1.38%  '

$ws.Range('D31').Value = '''6.82'
$ws.Range('E31').Value = '  +12.85%  '

$ws.Range('D32').Value = '''0.996'
$ws.Range('E32').Value = '  -0.38%  '

$ws.Range('D33').Value = '''19.11'
$ws.Range('E33').Value = '  +0.49%  '

$ws.Range('D34').Value = '''1.59'
$ws.Range('E34').Value = '  -0.26%  '

$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').Value = '''149.44'
$ws.Range('E35').Value = '  +0.13%  '

$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').Value = '''1.06'
$ws.Range('E36').Value = '  +19.20%  '

$ws.Range('E37').Value = '  +2.15%  '

$ws.Range('D38').Value = '''1.17'
$ws.Range('E38').Value = '  +0.64%  '

$ws.Range('D39').Value = '''0.879'
$ws.Range('E39').Value = '  +0.05%  '

$ws.Range('D40').Value = '''36.73'
$ws.Range('E40').Value = '  +0.47%  '

$ws.Range('E41').Value = '  +2.01%  '

$ws.Range('E42').Value = '  +0.48%  '

$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '''0.633'
$ws.Range('E43').Value = '  +0.83%  '

$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').Value = '''285.38'
$ws.Range('E44').Value = '  -1.86%  '

$ws.Range('E45').Value = '  -0.09%  '

$ws.Range('E46').Value = '  -0.57%  '

$ws.Range('D47').Value = '''19.92'
$ws.Range('E47').Value = '  +1.89%  '

$ws.Range('D48').Value = '''0.0548'
$ws.Range('E48').Value = '  -0.65%  '

$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''4.84'
$ws.Range('E49').Value = '  +1.88%  '

$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '''0.0234'
$ws.Range('E50').Value = '  -0.28%  '

$ws.Range('E51').Value = '  -0.99%  '
